# Refresh market-price snapshot values (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the Hades_Profits leve-crafting sheets, per the scheduled market-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 608
$ws.Range("I98").Value = 601.381
$ws.Range("J98").Value = 635.8
$ws.Range("K98").Value = 601.381
$ws.Range("L98").Value = 635.8
$ws.Range("M98").Value = 896.619
$ws.Range("N98").Value = -3631.8

$ws.Range("H116").Value = 2057.1428
$ws.Range("I116").Value = 1750
$ws.Range("J116").Value = 2466.6667
$ws.Range("K116").Value = 1750
$ws.Range("L116").Value = 2466.6667
$ws.Range("M116").Value = 1692
$ws.Range("N116").Value = -9350.6667

$ws.Range("H122").Value = 608
$ws.Range("I122").Value = 601.381
$ws.Range("J122").Value = 635.8
$ws.Range("K122").Value = 1804.143
$ws.Range("L122").Value = 1907.4
$ws.Range("M122").Value = 645.857
$ws.Range("N122").Value = -6807.4

$ws.Range("H129").Value = 914.08
$ws.Range("I129").Value = 544.7778
$ws.Range("J129").Value = 995.14636
$ws.Range("K129").Value = 1634.3334
$ws.Range("L129").Value = 2985.43908
$ws.Range("M129").Value = 3365.6666
$ws.Range("N129").Value = -12985.43908

$ws.Range("H132").Value = 1168217.9
$ws.Range("I132").Value = 1510.7949
$ws.Range("J132").Value = 16335410
$ws.Range("K132").Value = 4532.384700000001
$ws.Range("L132").Value = 49006230
$ws.Range("M132").Value = -2002.384700000001

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3327.3
$ws.Range("I63").Value = 2912.3333
$ws.Range("J63").Value = 3949.75
$ws.Range("K63").Value = 2912.3333
$ws.Range("L63").Value = 3949.75
$ws.Range("M63").Value = -2226.3333

$ws.Range("H66").Value = 3327.3
$ws.Range("I66").Value = 2912.3333
$ws.Range("J66").Value = 3949.75
$ws.Range("K66").Value = 14561.6665
$ws.Range("L66").Value = 19748.75
$ws.Range("M66").Value = -11129.6665

$ws.Range("H74").Value = 4201486
$ws.Range("I74").Value = 4922706
$ws.Range("J74").Value = 114572.22
$ws.Range("K74").Value = 4922706
$ws.Range("L74").Value = 114572.22
$ws.Range("M74").Value = -4921832

$ws.Range("H77").Value = 4201486
$ws.Range("I77").Value = 4922706
$ws.Range("J77").Value = 114572.22
$ws.Range("K77").Value = 24613530
$ws.Range("L77").Value = 572861.1
$ws.Range("M77").Value = -24609162

$ws.Range("H88").Value = 6806.6113
$ws.Range("I88").Value = 3585.5
$ws.Range("J88").Value = 9383.5
$ws.Range("K88").Value = 3585.5
$ws.Range("L88").Value = 9383.5
$ws.Range("M88").Value = -3179.5

$ws.Range("H91").Value = 6806.6113
$ws.Range("I91").Value = 3585.5
$ws.Range("J91").Value = 9383.5
$ws.Range("K91").Value = 3585.5
$ws.Range("L91").Value = 9383.5
$ws.Range("M91").Value = -2181.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11481.361
$ws.Range("I86").Value = 11176.32
$ws.Range("J86").Value = 12174.637
$ws.Range("K86").Value = 11176.32
$ws.Range("L86").Value = 12174.637
$ws.Range("M86").Value = -10053.32
$ws.Range("N86").Value = -14420.637

$ws.Range("H89").Value = 11481.361
$ws.Range("I89").Value = 11176.32
$ws.Range("J89").Value = 12174.637
$ws.Range("K89").Value = 55881.6
$ws.Range("L89").Value = 60873.185
$ws.Range("M89").Value = -50265.6
$ws.Range("N89").Value = -72105.185

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H107").Value = 360.64517
$ws.Range("I107").Value = 335.7143
$ws.Range("J107").Value = 413
$ws.Range("K107").Value = 335.7143
$ws.Range("L107").Value = 413
$ws.Range("M107").Value = 1584.2857
$ws.Range("N107").Value = -4253

$ws.Range("H134").Value = 33940.91
$ws.Range("I134").Value = 2420.926
$ws.Range("J134").Value = 155518
$ws.Range("K134").Value = 7262.778
$ws.Range("L134").Value = 466554
$ws.Range("M134").Value = -4727.778
$ws.Range("N134").Value = -471624

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4288059.5
$ws.Range("I4").Value = 30000000
$ws.Range("J4").Value = 2736
$ws.Range("K4").Value = 90000000
$ws.Range("L4").Value = 8208
$ws.Range("M4").Value = -89999888
$ws.Range("N4").Value = -8432

$ws.Range("H131").Value = 963.0909
$ws.Range("I131").Value = 462.45456
$ws.Range("J131").Value = 1034.6104
$ws.Range("K131").Value = 1387.36368
$ws.Range("L131").Value = 3103.8312
$ws.Range("M131").Value = 3652.63632
$ws.Range("N131").Value = -13183.8312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1260
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 1325
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 1325
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -4569

$ws.Range("H107").Value = 103
$ws.Range("I107").Value = 89.14286
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 89.14286
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 1830.85714
$ws.Range("N107").Value = -4040

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 9000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 88979.914
$ws.Range("I132").Value = 73017.14
$ws.Range("J132").Value = 113810.89
$ws.Range("K132").Value = 219051.42
$ws.Range("L132").Value = 341432.67
$ws.Range("M132").Value = -216521.42
$ws.Range("N132").Value = -346492.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2985.8948
$ws.Range("I16").Value = 1083.25
$ws.Range("J16").Value = 13133.333
$ws.Range("K16").Value = 1083.25
$ws.Range("L16").Value = 13133.333
$ws.Range("M16").Value = -913.25
$ws.Range("N16").Value = -13473.333

$ws.Range("H141").Value = 42357.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 42357.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 42357.5
$ws.Range("N141").Value = -52717.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 982.087
$ws.Range("I126").Value = 842.8570999999999
$ws.Range("J126").Value = 1198.6666
$ws.Range("K126").Value = 2528.5713
$ws.Range("L126").Value = 3595.9998
$ws.Range("M126").Value = -58.57129999999961

$ws.Range("H132").Value = 54113.316
$ws.Range("I132").Value = 39511.848
$ws.Range("J132").Value = 85749.836
$ws.Range("K132").Value = 118535.544
$ws.Range("L132").Value = 257249.508
$ws.Range("M132").Value = -116005.544
$ws.Range("N132").Value = -262309.508

$ws.Range("H136").Value = 43588.062
$ws.Range("I136").Value = 25949.85
$ws.Range("J136").Value = 144377.86
$ws.Range("K136").Value = 77849.54999999999
$ws.Range("L136").Value = 433133.58
$ws.Range("M136").Value = -75299.54999999999
$ws.Range("N136").Value = -438233.58

$ws.Range("H137").Value = 50800
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50800
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50800
$ws.Range("N137").Value = -61000

$ws.Range("H140").Value = 54363.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54363.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54363.332
$ws.Range("N140").Value = -64723.332

$ws.Range("H141").Value = 50083.332
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 50083.332
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 50083.332
$ws.Range("N141").Value = -60443.332
